# Add a new slide (slide 5) that mirrors the existing "Title and Content"
# slides (2-4): same title, single bold "Applications - ..." heading line.
#
# Duplicating slide 4 (the last "Title and Content" slide) and then
# overwriting the content placeholder's text is the simplest way to get
# PowerPoint to reuse the same paragraph/run formatting (pPr with
# lvl/indent/marL/spcBef/buNone, bold rPr) that the other slides already
# carry, instead of hand-building that XML.

$p = $ppt.ActivePresentation

$lastSlide = $p.Slides.Item($p.Slides.Count)
$newRange = $lastSlide.Duplicate()
$newSlide = $newRange.Item(1)

# Title placeholder keeps the same text as the other slides.
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Large Language Models on Graphs: A Comprehensive Survey - TKDE, December 2024"

# Content placeholder: replace the whole body with the single new heading.
# Assigning .Text collapses the body down to one paragraph while keeping
# the first paragraph's formatting (bold, spacing before, no bullet).
$newSlide.Shapes.Item(2).TextFrame.TextRange.Text = "Applications - Text-attributed graphs"
